{"js": "// Remove the standalone italic \"Ezekiel\" paragraph that immediately follows\n// the \"EZK\" Heading2 paragraph (the duplicate book-name subtitle line).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  const current = paragraphs.items[i];\n  if (current.text !== \"EZK\") {\n    continue;\n  }\n\n  const next = paragraphs.items[i + 1];\n  const nextRange = next.getRange();\n  nextRange.load(\"text\");\n  nextRange.font.load(\"italic\");\n  await context.sync();\n\n  if (nextRange.text === \"Ezekiel\" && nextRange.font.italic === true) {\n    target = next;\n    break;\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the standalone italic \"Ezekiel\" paragraph that immediately follows\n# the \"EZK\" Heading2 paragraph (the duplicate book-name subtitle line).\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -lt $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($ptext -ne \"EZK\") {\n        continue\n    }\n\n    $next = $d.Paragraphs.Item($i + 1)\n    $ntext = $next.Range.Text.TrimEnd([char]13, [char]7)\n    if ($ntext -eq \"Ezekiel\" -and $next.Range.Font.Italic) {\n        $next.Range.Delete()\n    }\n    break\n}\n"}
